$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Diebold-Mariano test statistic (C) and p-value (D) columns with corrected values
$ws.Range("C2").Value = -1.436217425223833
$ws.Range("D2").Value = 0.1650088613005696

$ws.Range("C3").Value = -0.9610255142477673
$ws.Range("D3").Value = 0.346983647897313

$ws.Range("C4").Value = -1.108421809652071
$ws.Range("D4").Value = 0.2796506844722986

$ws.Range("C5").Value = 0.5081432632920728
$ws.Range("D5").Value = 0.6164102838262626

$ws.Range("C6").Value = 0.5300394964552408
$ws.Range("D6").Value = 0.6013937355079939

$ws.Range("C7").Value = 0.5684900363574008
$ws.Range("D7").Value = 0.5754598973852505

$ws.Range("C8").Value = 2.011693928322679
$ws.Range("D8").Value = 0.05665866460692293
$ws.Range("G8").Value = "No"

$ws.Range("C9").Value = -0.1314882689649431
$ws.Range("D9").Value = 0.8965841351408768

$ws.Range("C10").Value = 1.306522625511748
$ws.Range("D10").Value = 0.2048692688869607

$ws.Range("C11").Value = 1.477219956278157
$ws.Range("D11").Value = 0.1537898730252496
